$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Widen column A from 11 to 22 (stored OOXML width).
# The COM layer's ColumnWidth setter adds a fixed 5/6 (0.8333...)
# character padding when serializing to the <col width="..."/> value,
# so compensate here to land on the exact target width of 22.
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668

# Update car names to include model names
$ws.Range("A2").Value = "Toyota Yaris"
$ws.Range("A3").Value = "Mazda MX-30"
$ws.Range("A4").Value = "Honda JAZZ"
$ws.Range("A5").Value = "Land Rover Defender"
$ws.Range("A6").Value = "SEAT Leon"
$ws.Range("A7").Value = "KIA Sorento"
$ws.Range("A8").Value = "Honda e"
$ws.Range("A9").Value = "Hyundai i10"
$ws.Range("A10").Value = "ISUZU D-Max Crew Cab"
$ws.Range("A11").Value = "Audi A3"
